$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Date conducted:" paragraph - trim the filler-space run by 3 chars
#    and append the conducted date/time as its own run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Date conducted: ")
$rng.Collapse(0)
$rng.InsertAfter("Nov 29, 2024, 7:27 PM")

# The run of 40 plain spaces that precedes "Date conducted: " needs to
# shrink to 37 spaces. Locate it via its known fixed offsets (plain
# Find on whitespace is unreliable because Word's Find treats runs of
# spaces/nbsp loosely) - it sits right before "Date conducted: ".
$dateRng = $d.Content
$dateRng.Find.Execute("Date conducted: ")
$spaceEnd = $dateRng.Start
$spaceStart = $spaceEnd - 40
$spaceRng = $d.Range($spaceStart, $spaceEnd)
$spaceRng.Text = "                                     "

# ---------------------------------------------------------------------
# 2. Test-case table: fill in Actual Output / Pass-Fail columns.
# ---------------------------------------------------------------------
$t = $d.Tables(1)

# Row 2 (Check if Dashboard is displayed upon login) already has a
# non-breaking-space placeholder in both cells; the new text is added
# alongside that placeholder (after it in Actual Output, before it in
# Pass/Fail), matching the original edit.
$c5 = $t.Cell(2, 5).Range
$insAfter = $d.Range($c5.End - 1, $c5.End - 1)
$insAfter.InsertAfter("Dashboard screen is displayed immediately after successful login")

$c6 = $t.Cell(2, 6).Range
$insBefore = $d.Range($c6.Start, $c6.Start)
$insBefore.InsertAfter("Pass")

# Rows 3-7 have genuinely empty cells in Actual Output / Pass-Fail.
$t.Cell(3, 5).Range.Text = "Dashboard is accompanied with all the UI components"
$t.Cell(3, 6).Range.Text = "Pass"

$t.Cell(4, 5).Range.Text = "Dashboard UI components retain its consistency with the changing window size"
$t.Cell(4, 6).Range.Text = "Pass"

$t.Cell(5, 5).Range.Text = "App shows settings, help and logout dialog buttons as options"
$t.Cell(5, 6).Range.Text = "Pass"

$t.Cell(6, 5).Range.Text = "Navigation arrows are displayed on each side while the main dashboard buttons are hidden"
$t.Cell(6, 6).Range.Text = "Pass"

$t.Cell(7, 5).Range.Text = "The displayed timer in the Dashboard starts countdown"
$t.Cell(7, 6).Range.Text = "Pass"

$t.Cell(8, 5).Range.Text = "There is no statistics screen that appears on Dashboard"
$t.Cell(8, 6).Range.Text = "Fail"

# ---------------------------------------------------------------------
# 3. Table column widths (tblGrid) resize.
# ---------------------------------------------------------------------
$t.Columns(1).Width = 164.45
$t.Columns(2).Width = 100.1
$t.Columns(3).Width = 54.1
$t.Columns(4).Width = 256.45
$t.Columns(5).Width = 237.4
